# The sheet contains a daily price log for "Agrícola del Norte S.A. de Arica - Mango".
# This commit adds one new daily record at the top of the data block (row 66),
# pushing every existing record below it down by one row (old row 144 becomes
# the new row 145). The new record itself is a fresh "Mango" price quote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 66 - shifts rows 66:144 down to 67:145.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new observation.
$ws.Range("A66").Value = 1
$ws.Range("B66").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C66").Value = "Arica y Parinacota"
$ws.Range("D66").Value = 44799
$ws.Range("E66").Value = 15
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100108
$ws.Range("H66").Value = "Tropicales y subtropicales"
$ws.Range("I66").Value = 100108002
$ws.Range("J66").Value = "Mango"
$ws.Range("K66").Value = "Sin especificar"
$ws.Range("L66").Value = "Especial"
$ws.Range("M66").Value = 456
$ws.Range("N66").Value = 9000
$ws.Range("O66").Value = 10000
$ws.Range("P66").Value = 9500
$ws.Range("Q66").Value = "$/bandeja 4 kilos"
$ws.Range("R66").Value = "Brasil"
$ws.Range("S66").Value = 2375
$ws.Range("T66").Value = 4

# Keep the date column's number format consistent with the rest of column D
# by copying the format from the row below (already shifted down).
$ws.Range("D66").NumberFormat = $ws.Range("D67").NumberFormat()
